# Updated symbol list on Wed Feb  1 18:49:55 UTC 2023 with GitHub Actions
# Refresh of coinranking.com snapshot: prices/volumes re-pulled and the
# coin roster (rows 6-18) rotated by one position, dropping GateToken off
# the front and appending it at the end with fresh data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.47%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.10%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.110"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.66%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07698"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.73%"

# Row 6
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.290"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.38%"

# Row 7
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.846"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.00%"

# Row 8
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.945"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.37%"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9207"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.05%"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1097"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-10.01%"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1843"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.11%"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08764"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.15%"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03345"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.87%"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09530"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.82%"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001378"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.41%"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006198"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "6.58%"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.362"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.46%"

# Row 18
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.384"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.80%"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3435"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.53%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.334"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "20.33%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1290"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.85%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2308"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-10.88%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04318"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.09%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004259"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.21%"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001330"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.01%"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002901"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02085"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.54%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04959"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.09%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007507"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.00%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1350"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.79%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008400"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.00%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002070"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.70%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008396"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.72%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006304"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-5.96%"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.01%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002851"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-14.26%"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001690"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.01%"
